$wb = $excel.ActiveWorkbook

$wsSPN = $wb.Worksheets.Item("SPN")
$wsITI = $wb.Worksheets.Item("ITI")

# --- Apply the "Semana" column formatting (style used in ITI!C2) to SPN's C2:C29 ---
$wsITI.Range("C2").Copy()
$wsSPN.Range("C2:C29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Set column C (Semana) to numeric value 2 for data rows in both sheets ---
$wsSPN.Range("C2:C29").Value = 2
$wsITI.Range("C2:C22").Value = 2

# --- Selection / active sheet changes ---
# ITI sheet selection becomes C2:C22, active cell defaults to top-left (C2)
$wsITI.Activate()
$wsITI.Range("C2:C22").Select()

# SPN becomes the active sheet with C2:C29 selected, active cell defaults to top-left (C2)
$wsSPN.Activate()
$wsSPN.Range("C2:C29").Select()
